$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.340.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.15%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.055.45"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.49%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.87%  "

# Row 6
$ws.Range("E6").Value = "  +2.73%  "

# Row 7
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.81"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.36%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  +3.55%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.69"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.55%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0761"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.60%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.102"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.356.42"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.36"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.28%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.86"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.774"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.12%  "

# Row 17
$ws.Range("E17").Value = "  +1.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.041.67"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.82%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.505.82"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +15.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.19"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0811"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.66%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "225.63"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.37%  "

# Row 24
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("E25").Value = "  +0.30%  "

# Row 26
$ws.Range("E26").Value = "  +0.66%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.27"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.12%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.48"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.13%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.90"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.36%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.130"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.63%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.09"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.87%  "

# Row 32
$ws.Range("E32").Value = "  +0.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.48"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.54%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0621"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.72%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.57"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.58%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.59"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.18%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.06%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.93"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.97%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.31"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.62%  "

# Row 40
$ws.Range("E40").Value = "  -0.93%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.68"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +11.12%  "

# Row 42
$ws.Range("E42").Value = "  -0.87%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0945"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.71%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.83"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.27%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.456.39"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.04%  "

# Row 46
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.43%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0212"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.65%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.66"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.18%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.03"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.43%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.17"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.94"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.11%  "
